# ExcelData provider issue fixed.
# Remove the "Login" and "Sheet2" worksheets, keep the third sheet
# (originally "Sheet3"), rename it to "Login", move it to the first
# position, and populate it with the new test data.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Login").Delete()
$wb.Worksheets.Item("Sheet2").Delete()

$ws = $wb.Worksheets.Item("Sheet3")
$ws.Name = "Login"
$ws.Move($wb.Worksheets.Item(1))
$ws.Activate()
$ws.Select()

$ws.Range("A1").Value = "abcdefg"
$ws.Range("B1").Value = "dfdsdsd"
$ws.Range("B1").Select()
